$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price column D, Volume(1h) column E) per the
# 2024-06-27 16:49:17 UTC GitHub Actions refresh.

$ws.Range("D2").Value = "61.759.42"
$ws.Range("E2").Value = "  +0.51%  "

$ws.Range("D3").Value = "3.451.95"
$ws.Range("E3").Value = "  +2.74%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.83"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.91"
$ws.Range("E6").Value = "  +8.48%  "

$ws.Range("D7").Value = "3.452.18"
$ws.Range("E7").Value = "  +2.79%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.68"
$ws.Range("E10").Value = "  +2.89%  "

$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("E12").Value = "  -0.25%  "

$ws.Range("D13").Value = "4.037.55"
$ws.Range("E13").Value = "  +2.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.95"
$ws.Range("E14").Value = "  +8.20%  "

$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").Value = "3.446.63"
$ws.Range("E17").Value = "  +2.53%  "

$ws.Range("D18").Value = "61.805.00"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.29"
$ws.Range("E19").Value = "  +7.74%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.10"
$ws.Range("E20").Value = "  +0.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.39"
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "382.50"
$ws.Range("E22").Value = "  +0.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.566"
$ws.Range("E23").Value = "  +2.43%  "

$ws.Range("D24").Value = "3.589.78"
$ws.Range("E24").Value = "  +2.54%  "

$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.28"
$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("E28").Value = "  -0.67%  "

$ws.Range("E29").Value = "  +8.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("E30").Value = "  +2.95%  "

$ws.Range("E31").Value = "  -12.49%  "

$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.20"
$ws.Range("E33").Value = "  +0.26%  "

$ws.Range("E34").Value = "  +1.07%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.02"
$ws.Range("E36").Value = "  +2.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.02"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("E38").Value = "  -0.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("E39").Value = "  +2.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "165.86"
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0783"
$ws.Range("E41").Value = "  +2.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.792"
$ws.Range("E42").Value = "  +3.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.84"
$ws.Range("E43").Value = "  +7.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("E45").Value = "  +0.77%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.16"
$ws.Range("E46").Value = "  +1.66%  "

$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").Value = "2.610.72"
$ws.Range("E49").Value = "  +10.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.54"
$ws.Range("E50").Value = "  +1.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.86"
$ws.Range("E51").Value = "  +0.54%  "

